$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price (D) and volume (E) columns remain plain text so that
# numeric-looking strings (e.g. "0.660", "0.0000272") are not coerced
# into numbers/scientific notation by Excel when the Value is assigned.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '67.884.98'
$ws.Range('E2').Value = '  -0.60%  '
$ws.Range('D3').Value = '3.531.73'
$ws.Range('E3').Value = '  -2.99%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').Value = '199.52'
$ws.Range('E5').Value = '  -0.96%  '
$ws.Range('D6').Value = '556.42'
$ws.Range('E6').Value = '  -2.85%  '
$ws.Range('E7').Value = '  +2.69%  '
$ws.Range('D8').Value = '3.524.78'
$ws.Range('E8').Value = '  -3.04%  '
$ws.Range('D10').Value = '0.660'
$ws.Range('E10').Value = '  -2.51%  '
$ws.Range('D11').Value = '62.27'
$ws.Range('E11').Value = '  +8.29%  '
$ws.Range('D12').Value = '0.145'
$ws.Range('E12').Value = '  -6.01%  '
$ws.Range('D13').Value = '0.0000272'
$ws.Range('E13').Value = '  -7.29%  '
$ws.Range('D14').Value = '9.98'
$ws.Range('E14').Value = '  -0.96%  '
$ws.Range('D15').Value = '4.086.21'
$ws.Range('E15').Value = '  -3.20%  '
$ws.Range('D16').Value = '3.524.43'
$ws.Range('E16').Value = '  -3.29%  '
$ws.Range('E17').Value = '  -1.71%  '
$ws.Range('D18').Value = '67.608.77'
$ws.Range('E18').Value = '  -0.93%  '
$ws.Range('D19').Value = '18.43'
$ws.Range('E19').Value = '  -0.99%  '
$ws.Range('D20').Value = '11.90'
$ws.Range('E20').Value = '  -4.71%  '
$ws.Range('D21').Value = '1.03'
$ws.Range('E21').Value = '  -5.05%  '
$ws.Range('D22').Value = '396.81'
$ws.Range('E22').Value = '  -1.46%  '
$ws.Range('D23').Value = '4.01'
$ws.Range('E23').Value = '  -4.93%  '
$ws.Range('D24').Value = '11.89'
$ws.Range('E24').Value = '  -8.91%  '
$ws.Range('D25').Value = '85.29'
$ws.Range('E25').Value = '  -0.89%  '
$ws.Range('B26').Value = 'InternetComputer(DFINITY)'
$ws.Range('C26').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D26').Value = '12.44'
$ws.Range('E26').Value = '  -1.48%  '
$ws.Range('B27').Value = 'Toncoin'
$ws.Range('C27').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D27').Value = '3.88'
$ws.Range('E27').Value = '  +0.87%  '
$ws.Range('D28').Value = '2.84'
$ws.Range('E28').Value = '  -4.34%  '
$ws.Range('D29').Value = '8.92'
$ws.Range('E29').Value = '  -2.61%  '
$ws.Range('D30').Value = '725.15'
$ws.Range('E30').Value = '  +3.90%  '
$ws.Range('D31').Value = '31.37'
$ws.Range('E31').Value = '  -1.90%  '
$ws.Range('D32').Value = '7.09'
$ws.Range('E32').Value = '  -14.02%  '
$ws.Range('D33').Value = '11.77'
$ws.Range('E33').Value = '  -3.89%  '
$ws.Range('D34').Value = '64.52'
$ws.Range('E34').Value = '  +0.12%  '
$ws.Range('D35').Value = '0.112'
$ws.Range('E35').Value = '  -4.07%  '
$ws.Range('D36').Value = '38.80'
$ws.Range('E36').Value = '  -9.41%  '
$ws.Range('E37').Value = '  -0.01%  '
$ws.Range('D38').Value = '0.396'
$ws.Range('E38').Value = '  -7.48%  '
$ws.Range('D39').Value = '0.132'
$ws.Range('E39').Value = '  -6.16%  '
$ws.Range('D40').Value = '3.03'
$ws.Range('E40').Value = '  -3.62%  '
$ws.Range('B41').Value = 'Maker'
$ws.Range('C41').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D41').Value = '3.084.54'
$ws.Range('E41').Value = '  -5.15%  '
$ws.Range('B42').Value = 'FirstDigitalUSD'
$ws.Range('C42').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D42').Value = '0.997'
$ws.Range('E42').Value = '  -0.15%  '
$ws.Range('D43').Value = '0.0₃0687'
$ws.Range('E43').Value = '  -13.04%  '
$ws.Range('B44').Value = 'Fetch.AI'
$ws.Range('C44').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D44').Value = '2.51'
$ws.Range('E44').Value = '  -10.45%  '
$ws.Range('B45').Value = 'WEMIXToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D45').Value = '2.77'
$ws.Range('E45').Value = '  +2.66%  '
$ws.Range('D46').Value = '0.0409'
$ws.Range('E46').Value = '  -2.78%  '
$ws.Range('D47').Value = '0.132'
$ws.Range('E47').Value = '  +0.09%  '
$ws.Range('D48').Value = '139.28'
$ws.Range('E48').Value = '  -2.04%  '
$ws.Range('D49').Value = '2.57'
$ws.Range('E49').Value = '  -13.87%  '
$ws.Range('D50').Value = '2.95'
$ws.Range('E50').Value = '  -3.18%  '
$ws.Range('D51').Value = '8.26'
$ws.Range('E51').Value = '  -7.58%  '

Write-Output "Applied cryptos list update"
